$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Dorados de Sinaloa - Tlaxcala) odds update
$ws.Range("I6").Value = 3.15
$ws.Range("J6").Value = 2.72
$ws.Range("L6").Value = 3.65
$ws.Range("Q6").Value = 1.95
$ws.Range("U6").Value = 1.7
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 7.8
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 8.75
$ws.Range("AA6").Value = 17.5
$ws.Range("AB6").Value = 27
$ws.Range("AG6").Value = 9
$ws.Range("AH6").Value = 16
$ws.Range("AK6").Value = 29
$ws.Range("AL6").Value = 37
$ws.Range("AN6").Value = 4.15
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 17.5
$ws.Range("AQ6").Value = 40
$ws.Range("AR6").Value = 65
$ws.Range("AU6").Value = 6.6
$ws.Range("AW6").Value = 5.1
$ws.Range("AX6").Value = 17
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 80
$ws.Range("BA6").Value = 110
$ws.Range("BB6").Value = 300

# Row 10 (Sportivo Trinidense - Libertad Asuncion) odds update
$ws.Range("G10").Value = 3.7
$ws.Range("H10").Value = 3.25
$ws.Range("I10").Value = 1.91
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.75
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("Y10").Value = 13
$ws.Range("AC10").Value = 8
$ws.Range("AE10").Value = 17
$ws.Range("AG10").Value = 6.5
$ws.Range("AH10").Value = 8.5
$ws.Range("AJ10").Value = 17
$ws.Range("AN10").Value = 5.5
$ws.Range("AR10").Value = 101
$ws.Range("AT10").Value = 2.5
$ws.Range("AW10").Value = 4
$ws.Range("BA10").Value = 67

# Row 12 (Academico Viseu - Maritimo) odds update
$ws.Range("N12").Value = 9.5
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
$ws.Range("Q12").Value = 2.05
$ws.Range("R12").Value = 1.75
$ws.Range("S12").Value = 1.44
$ws.Range("T12").Value = 2.63
$ws.Range("AT12").Value = 2.63

# Row 26 (Carabobo - Dep. Tachira) odds update
$ws.Range("G26").Value = 2.77
$ws.Range("I26").Value = 2.77
$ws.Range("J26").Value = 3.45
$ws.Range("L26").Value = 3.5
$ws.Range("N26").Value = 5.8
$ws.Range("O26").Value = 1.52
$ws.Range("P26").Value = 2.22
$ws.Range("Q26").Value = 2.47
$ws.Range("U26").Value = 2
$ws.Range("W26").Value = 6.6
$ws.Range("X26").Value = 12.5
$ws.Range("Y26").Value = 10.75
$ws.Range("Z26").Value = 35
$ws.Range("AA26").Value = 29
$ws.Range("AB26").Value = 45
$ws.Range("AD26").Value = 5.5
$ws.Range("AG26").Value = 6.5
$ws.Range("AH26").Value = 12.5
$ws.Range("AI26").Value = 10.75
$ws.Range("AJ26").Value = 35
$ws.Range("AK26").Value = 30
$ws.Range("AL26").Value = 50
$ws.Range("AN26").Value = 4.45
$ws.Range("AO26").Value = 16
$ws.Range("AP26").Value = 26
$ws.Range("AQ26").Value = 80
$ws.Range("AS26").Value = 400
$ws.Range("AW26").Value = 4.45
$ws.Range("AX26").Value = 16
$ws.Range("AY26").Value = 27
$ws.Range("AZ26").Value = 80
$ws.Range("BA26").Value = 150
$ws.Range("BB26").Value = 450
